# Generate Report for Handoff
#
# The workbook tracks localization handoff status for a repo. A new file
# (dbc24d4f-a94b-48c6-96b2-25dae80136a8.md, plus two screenshots referenced
# from it) has been sent out for localization, so each sheet grows two new
# data rows (inserted just above the final ".localization-config" row) and
# the existing first data row's source file is swapped for a different one
# that's now also "Ready for handoff" / has a dependency.

$wb = $excel.ActiveWorkbook

function Clear-AllHyperlinks($ws) {
    if ($ws.Hyperlinks.Count -gt 0) {
        $ws.Range("A1").Hyperlinks.Delete()
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview" (3 columns: File Name / zh-cn / de-de)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

Clear-AllHyperlinks $ws1

# Make room for two new file rows above the ".localization-config" row.
$ws1.Rows("3:4").Insert()

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/1ff5ce0f-ff33-49cb-8c3a-51dcf2350895.png", "", "", "1ff5ce0f-ff33-49cb-8c3a-51dcf2350895.png") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/c2f4a5ca-e29d-4279-83d1-5951ea224f02.png", "", "", "c2f4a5ca-e29d-4279-83d1-5951ea224f02.png") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/dbc24d4f-a94b-48c6-96b2-25dae80136a8.md", "", "", "dbc24d4f-a94b-48c6-96b2-25dae80136a8.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheets 2 & 3: "zh-cn" / "de-de" (9-column per-language detail sheets)
# ---------------------------------------------------------------------
$langs = @(
    @{
        SheetName  = "zh-cn"
        C2         = "588848e5b175a00a074e02b26f9fcd22491686b4.png"
        D2         = "2016-03-08 14:53:13"
        C3         = "c8cac34cddcea9cdf96896c0d52765a4abdb3fcb.png"
        D3         = "2016-03-08 14:53:13"
        C4         = "dbc24d4f-a94b-48c6-96b2-25dae80136a8.541924b48566044bd45bb80c4e737f7553c213b5.zh-cn.xlf"
        D4         = "2016-03-08 14:53:13"
        HtPath     = "oltest.zh-cn"
    },
    @{
        SheetName  = "de-de"
        C2         = "588848e5b175a00a074e02b26f9fcd22491686b4.png"
        D2         = "2016-03-08 14:53:19"
        C3         = "c8cac34cddcea9cdf96896c0d52765a4abdb3fcb.png"
        D3         = "2016-03-08 14:53:19"
        C4         = "dbc24d4f-a94b-48c6-96b2-25dae80136a8.541924b48566044bd45bb80c4e737f7553c213b5.de-de.xlf"
        D4         = "2016-03-08 14:53:19"
        HtPath     = "oltest.de-de"
    }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.SheetName)

    Clear-AllHyperlinks $ws

    # Make room for two new file rows above the ".localization-config" row.
    $ws.Rows("3:4").Insert()

    # --- Row 2 (existing source file, now has a dependency) ---
    $ws.Range("D2").Value = $lang.D2
    $ws.Range("H2").Value = "IsDependency"
    $ws.Range("I2").Value = "e2e\dbc24d4f-a94b-48c6-96b2-25dae80136a8.md"

    # --- Row 3 (new: second screenshot, depends on the .md handoff) ---
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $lang.D3
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "IsDependency"
    $ws.Range("I3").Value = "e2e\dbc24d4f-a94b-48c6-96b2-25dae80136a8.md"

    # --- Row 4 (new: the markdown file itself, handed off normally) ---
    $ws.Range("B4").Value = "Ready for handoff"
    $ws.Range("D4").Value = $lang.D4
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Include"

    # --- Row 5 (was row 3: ".localization-config", unchanged content) ---
    # values already correct after the row insert/shift; only the
    # hyperlinks (cleared above) need to be re-created below.

    $htBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/$($lang.HtPath)/ci/ht"

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/1ff5ce0f-ff33-49cb-8c3a-51dcf2350895.png", "", "", "1ff5ce0f-ff33-49cb-8c3a-51dcf2350895.png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), "$htBase/$($lang.C2)", "", "", $lang.C2) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/c2f4a5ca-e29d-4279-83d1-5951ea224f02.png", "", "", "c2f4a5ca-e29d-4279-83d1-5951ea224f02.png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), "$htBase/$($lang.C3)", "", "", $lang.C3) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/dbc24d4f-a94b-48c6-96b2-25dae80136a8.md", "", "", "dbc24d4f-a94b-48c6-96b2-25dae80136a8.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C4"), "$htBase/$($lang.C4)", "", "", $lang.C4) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/.localization-config", "", "", ".localization-config") | Out-Null
}

Write-Host "Localization status report updated for handoff."
